$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39-45 down to 40-46.
$ws.Rows.Item(39).Insert()

# Populate the new row 39 with the new record (copy of row 40's static
# columns, but with the updated measurement values from the diff).
$ws.Cells.Item(39, 1).Value = 3
$ws.Cells.Item(39, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(39, 3).Value = "Coquimbo"
$ws.Cells.Item(39, 4).Value = 44543
$ws.Cells.Item(39, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(39, 5).Value = 5
$ws.Cells.Item(39, 6).Value = 100112022
$ws.Cells.Item(39, 7).Value = "Arveja Verde"
$ws.Cells.Item(39, 8).Value = "Perfection"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 35
$ws.Cells.Item(39, 11).Value = 18000
$ws.Cells.Item(39, 12).Value = 18000
$ws.Cells.Item(39, 13).Value = 18000
$ws.Cells.Item(39, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(39, 16).Value = 720
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"
